$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.899.76"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "2.301.07"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'305.68"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'97.12"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "'35.61"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D12").Value = "'18.28"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "'0.119"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "2.659.37"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "2.298.29"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "42.813.96"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").Value = "'67.56"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "'236.62"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").Value = "'2.46"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").Value = "'4.02"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'25.50"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").Value = "'167.32"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "'2.07"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "'33.01"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").Value = "'5.01"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("D36").Value = "'17.38"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").Value = "2.009.27"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "'18.44"
$ws.Range("E45").Value = "  +5.18%  "
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("E49").Value = "  +7.55%  "
$ws.Range("D50").Value = "'53.96"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "2.527.61"
$ws.Range("E51").Value = "  -0.42%  "
